$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (sex) was blank for data rows 2-6; fill it with "U" (unknown/unsexed).
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = "U"
}

# Columns F (sire) and G (dam) data values are no longer used for this cross table;
# clear the old numeric codes from rows 2-6 (headers in row 1 stay).
$ws.Range("F2:G6").ClearContents()

# The active selection moved to D7 after the edit.
$ws.Range("D7").Select()
